# Regenerate merged AHB files
# - rename header strings "_old" -> "_FV2410" and "_new" -> "_FV2504"
# - turn the used range A1:U57 into an Excel Table (Table1)
# - freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2410 = "_FV2410"
$fv2504 = "_FV2504"

# 1) Rename the header row (row 1, columns A:U) in place -- do this BEFORE
#    creating the table so the ListObject's column names pick up the new text.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith($oldSuffix)) {
            $cell.Value = $val.Substring(0, $val.Length - $oldSuffix.Length) + $fv2410
        } elseif ($val.EndsWith($newSuffix)) {
            $cell.Value = $val.Substring(0, $val.Length - $newSuffix.Length) + $fv2504
        }
    }
}

# 2) Convert the used range into an Excel Table.
$tableRange = $ws.Range("A1:U57")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# 3) Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
